# Update loading_percent values for case with 380 kV
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 12.10079634308037
$ws.Range("C2").Value2 = 7.806623520735326
$ws.Range("D2").Value2 = 5.532385992127314
$ws.Range("E2").Value2 = 11.58203427379016
$ws.Range("F2").Value2 = 26.63433572022027
$ws.Range("I2").Value2 = 23.91549173517083
$ws.Range("K2").Value2 = 9.43443779019503
$ws.Range("L2").Value2 = 9.766658145443696
$ws.Range("N2").Value2 = 19.15499729034599
$ws.Range("O2").Value2 = 23.90746244434957
$ws.Range("B3").Value2 = 11.84286945263442
$ws.Range("C3").Value2 = 7.771582587444184
$ws.Range("D3").Value2 = 5.479148507636935
$ws.Range("E3").Value2 = 11.59103829459997
$ws.Range("F3").Value2 = 26.66347591426661
$ws.Range("I3").Value2 = 23.99452504884313
$ws.Range("K3").Value2 = 9.252926683634682
$ws.Range("L3").Value2 = 9.751873418654027
$ws.Range("N3").Value2 = 19.21249427670219
$ws.Range("O3").Value2 = 23.97201195760561
$ws.Range("B4").Value2 = 11.68376696587798
$ws.Range("C4").Value2 = 7.749725850356129
$ws.Range("D4").Value2 = 5.445647036927015
$ws.Range("E4").Value2 = 11.59861685727043
$ws.Range("F4").Value2 = 26.68775224722922
$ws.Range("I4").Value2 = 24.04722864910934
$ws.Range("K4").Value2 = 9.141145780802068
$ws.Range("L4").Value2 = 9.744480418479293
$ws.Range("N4").Value2 = 19.249459744403
$ws.Range("O4").Value2 = 24.01630033891879
$ws.Range("B5").Value2 = 11.61884015142188
$ws.Range("C5").Value2 = 7.740735523568848
$ws.Range("D5").Value2 = 5.431795931979114
$ws.Range("E5").Value2 = 11.60222110865911
$ws.Range("F5").Value2 = 26.69924887635827
$ws.Range("I5").Value2 = 24.06975541000665
$ws.Range("K5").Value2 = 9.095571915532986
$ws.Range("L5").Value2 = 9.741893993363067
$ws.Range("N5").Value2 = 19.2649426252255
$ws.Range("O5").Value2 = 24.03551708753019
$ws.Range("B6").Value2 = 11.60805637304841
$ws.Range("C6").Value2 = 7.739237709746359
$ws.Range("D6").Value2 = 5.429484112745891
$ws.Range("E6").Value2 = 11.60285076176933
$ws.Range("F6").Value2 = 26.70125469272138
$ws.Range("I6").Value2 = 24.07355933556625
$ws.Range("K6").Value2 = 9.088004884535614
$ws.Range("L6").Value2 = 9.74149033632621
$ws.Range("N6").Value2 = 19.26753889712439
$ws.Range("O6").Value2 = 24.03877856229984
$ws.Range("B7").Value2 = 11.68289158428295
$ws.Range("C7").Value2 = 7.749604938791737
$ws.Range("D7").Value2 = 5.445461034353348
$ws.Range("E7").Value2 = 11.59866337605382
$ws.Range("F7").Value2 = 26.68790080346044
$ws.Range("I7").Value2 = 24.04752820394385
$ws.Range("K7").Value2 = 9.140531163766783
$ws.Range("L7").Value2 = 9.744443807813489
$ws.Range("N7").Value2 = 19.24966685321439
$ws.Range("O7").Value2 = 24.01655477218306
$ws.Range("B8").Value2 = 12.01207071404892
$ws.Range("C8").Value2 = 7.7946136757867
$ws.Range("D8").Value2 = 5.51420184035591
$ws.Range("E8").Value2 = 11.58471366770218
$ws.Range("F8").Value2 = 26.64305751224538
$ws.Range("I8").Value2 = 23.94187511864747
$ws.Range("K8").Value2 = 9.371958097377938
$ws.Range("L8").Value2 = 9.76121232099957
$ws.Range("N8").Value2 = 19.17447801559574
$ws.Range("O8").Value2 = 23.92875210476643
$ws.Range("B9").Value2 = 12.6478160472391
$ws.Range("C9").Value2 = 7.880065451616209
$ws.Range("D9").Value2 = 5.642310592605126
$ws.Range("E9").Value2 = 11.57360080537385
$ws.Range("F9").Value2 = 26.60581671573273
$ws.Range("I9").Value2 = 23.76786090017813
$ws.Range("K9").Value2 = 9.820533568007592
$ws.Range("L9").Value2 = 9.80733807796272
$ws.Range("N9").Value2 = 19.04016478107055
$ws.Range("O9").Value2 = 23.79357204852329
$ws.Range("B10").Value2 = 13.10391677660666
$ws.Range("C10").Value2 = 7.94100103724081
$ws.Range("D10").Value2 = 5.732044923334892
$ws.Range("E10").Value2 = 11.5752981343789
$ws.Range("F10").Value2 = 26.60937968931779
$ws.Range("I10").Value2 = 23.66027476616907
$ws.Range("K10").Value2 = 10.14353965717149
$ws.Range("L10").Value2 = 9.849114170829568
$ws.Range("N10").Value2 = 18.94941055842205
$ws.Range("O10").Value2 = 23.71689873241209
$ws.Range("B11").Value2 = 13.3080477727507
$ws.Range("C11").Value2 = 7.968293261224519
$ws.Range("D11").Value2 = 5.771852039404273
$ws.Range("E11").Value2 = 11.57819981610829
$ws.Range("F11").Value2 = 26.61770420217636
$ws.Range("I11").Value2 = 23.61573911751093
$ws.Range("K11").Value2 = 10.28839555327013
$ws.Range("L11").Value2 = 9.869786500057007
$ws.Range("N11").Value2 = 18.90982808009804
$ws.Range("O11").Value2 = 23.68695088983071
$ws.Range("B12").Value2 = 13.38479004859974
$ws.Range("C12").Value2 = 7.978564171255329
$ws.Range("D12").Value2 = 5.786775071971398
$ws.Range("E12").Value2 = 11.57960350616716
$ws.Range("F12").Value2 = 26.62181816394605
$ws.Range("I12").Value2 = 23.59950889257119
$ws.Range("K12").Value2 = 10.34289794469997
$ws.Range("L12").Value2 = 9.877850148114655
$ws.Range("N12").Value2 = 18.89508272776666
$ws.Range("O12").Value2 = 23.6763206027537
$ws.Range("B13").Value2 = 13.36828818258214
$ws.Range("C13").Value2 = 7.976355036953383
$ws.Range("D13").Value2 = 5.78356793132672
$ws.Range("E13").Value2 = 11.5792876597995
$ws.Range("F13").Value2 = 26.62088942460437
$ws.Range("I13").Value2 = 23.60297613003946
$ws.Range("K13").Value2 = 10.33117630523461
$ws.Range("L13").Value2 = 9.876103096783496
$ws.Range("N13").Value2 = 18.8982475829179
$ws.Range("O13").Value2 = 23.67857841651353
$ws.Range("B14").Value2 = 13.31437299903569
$ws.Range("C14").Value2 = 7.96913955374393
$ws.Range("D14").Value2 = 5.773082823634621
$ws.Range("E14").Value2 = 11.57830919649062
$ws.Range("F14").Value2 = 26.61802340357051
$ws.Range("I14").Value2 = 23.61439112654329
$ws.Range("K14").Value2 = 10.29288682794576
$ws.Range("L14").Value2 = 9.870445207883719
$ws.Range("N14").Value2 = 18.90861009367873
$ws.Range("O14").Value2 = 23.68606208768923
$ws.Range("B15").Value2 = 13.28127360128556
$ws.Range("C15").Value2 = 7.964711443115757
$ws.Range("D15").Value2 = 5.766640564892814
$ws.Range("E15").Value2 = 11.57774952172614
$ws.Range("F15").Value2 = 26.61639303535638
$ws.Range("I15").Value2 = 23.62146579319928
$ws.Range("K15").Value2 = 10.26938616094701
$ws.Range("L15").Value2 = 9.867010119700828
$ws.Range("N15").Value2 = 18.91498913221535
$ws.Range("O15").Value2 = 23.69073858887813
$ws.Range("B16").Value2 = 13.09050230327968
$ws.Range("C16").Value2 = 7.939208600958775
$ws.Range("D16").Value2 = 5.729422582481663
$ws.Range("E16").Value2 = 11.57515124425598
$ws.Range("F16").Value2 = 26.60897037680411
$ws.Range("I16").Value2 = 23.66327401007924
$ws.Range("K16").Value2 = 10.13402658029048
$ws.Range("L16").Value2 = 9.847796383836616
$ws.Range("N16").Value2 = 18.95203152371468
$ws.Range("O16").Value2 = 23.71895523006378
$ws.Range("B17").Value2 = 12.97255638267748
$ws.Range("C17").Value2 = 7.923452068694065
$ws.Range("D17").Value2 = 5.706327144515662
$ws.Range("E17").Value2 = 11.57410181903323
$ws.Range("F17").Value2 = 26.60613238297477
$ws.Range("I17").Value2 = 23.69005108738326
$ws.Range("K17").Value2 = 10.05041677820497
$ws.Range("L17").Value2 = 9.836433597326998
$ws.Range("N17").Value2 = 18.9751909988872
$ws.Range("O17").Value2 = 23.7375291406822
$ws.Range("B18").Value2 = 12.9044057612781
$ws.Range("C18").Value2 = 7.914349225858158
$ws.Range("D18").Value2 = 5.692948174549902
$ws.Range("E18").Value2 = 11.57369875327278
$ws.Range("F18").Value2 = 26.60513121720914
$ws.Range("I18").Value2 = 23.70586719831731
$ws.Range("K18").Value2 = 10.00213366866289
$ws.Range("L18").Value2 = 9.830055320117063
$ws.Range("N18").Value2 = 18.98867198763458
$ws.Range("O18").Value2 = 23.74867645158696
$ws.Range("B19").Value2 = 12.8812800660342
$ws.Range("C19").Value2 = 7.911260350286615
$ws.Range("D19").Value2 = 5.68840211892694
$ws.Range("E19").Value2 = 11.57359676227037
$ws.Range("F19").Value2 = 26.60490070254779
$ws.Range("I19").Value2 = 23.71129344923873
$ws.Range("K19").Value2 = 9.985754333433098
$ws.Range("L19").Value2 = 9.827922890472667
$ws.Range("N19").Value2 = 18.99326397869758
$ws.Range("O19").Value2 = 23.75253040790471
$ws.Range("B20").Value2 = 12.98514470506759
$ws.Range("C20").Value2 = 7.92513354747163
$ws.Range("D20").Value2 = 5.708795573271746
$ws.Range("E20").Value2 = 11.57419278657865
$ws.Range("F20").Value2 = 26.60636918260634
$ws.Range("I20").Value2 = 23.68715770100486
$ws.Range("K20").Value2 = 10.05933754245591
$ws.Range("L20").Value2 = 9.837626936683927
$ws.Range("N20").Value2 = 18.97270905362501
$ws.Range("O20").Value2 = 23.73550387549545
$ws.Range("B21").Value2 = 13.33022490981768
$ws.Range("C21").Value2 = 7.971260675105526
$ws.Range("D21").Value2 = 5.776166695241125
$ws.Range("E21").Value2 = 11.5785883319521
$ws.Range("F21").Value2 = 26.61883914703964
$ws.Range("I21").Value2 = 23.6110210392484
$ws.Range("K21").Value2 = 10.30414331252311
$ws.Range("L21").Value2 = 9.872100712351688
$ws.Range("N21").Value2 = 18.90555976843383
$ws.Range("O21").Value2 = 23.68384466560761
$ws.Range("B22").Value2 = 13.55247251133385
$ws.Range("C22").Value2 = 8.001032829154529
$ws.Range("D22").Value2 = 5.819314777366527
$ws.Range("E22").Value2 = 11.58323743435738
$ws.Range("F22").Value2 = 26.63259257871161
$ws.Range("I22").Value2 = 23.56495984128346
$ws.Range("K22").Value2 = 10.46206784295275
$ws.Range("L22").Value2 = 9.896001682267306
$ws.Range("N22").Value2 = 18.86309364759249
$ws.Range("O22").Value2 = 23.65422306656894
$ws.Range("B23").Value2 = 13.43417882414672
$ws.Range("C23").Value2 = 7.985177989911685
$ws.Range("D23").Value2 = 5.796368331222883
$ws.Range("E23").Value2 = 11.58059407401657
$ws.Range("F23").Value2 = 26.62474035947553
$ws.Range("I23").Value2 = 23.58920486294765
$ws.Range("K23").Value2 = 10.37798639862419
$ws.Range("L23").Value2 = 9.883121417262931
$ws.Range("N23").Value2 = 18.88562905697019
$ws.Range("O23").Value2 = 23.66965344220572
$ws.Range("B24").Value2 = 12.9794545837632
$ws.Range("C24").Value2 = 7.924373488414011
$ws.Range("D24").Value2 = 5.707679911113706
$ws.Range("E24").Value2 = 11.57415103626326
$ws.Range("F24").Value2 = 26.60626016156998
$ws.Range("I24").Value2 = 23.6884644877552
$ws.Range("K24").Value2 = 10.05530512974318
$ws.Range("L24").Value2 = 9.837086946882796
$ws.Range("N24").Value2 = 18.97383062297908
$ws.Range("O24").Value2 = 23.73641803731224
$ws.Range("B25").Value2 = 12.47741139755548
$ws.Range("C25").Value2 = 7.857263008752597
$ws.Range("D25").Value2 = 5.60840241062929
$ws.Range("E25").Value2 = 11.57487169976241
$ws.Range("F25").Value2 = 26.61045916878754
$ws.Range("I25").Value2 = 23.81138094591081
$ws.Range("K25").Value2 = 9.700100328178921
$ws.Range("L25").Value2 = 9.793460124805851
$ws.Range("N25").Value2 = 19.07510250535235
$ws.Range("O25").Value2 = 23.82617115558725
